$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.197657823562622
$ws.Range("B1").Value = 2.053754329681396
$ws.Range("C1").Value = 4.331416130065918
$ws.Range("D1").Value = 3.025156021118164
$ws.Range("E1").Value = 1.20698881149292
